# Update summary charts / summary reports per Prof. Erhardt's comments.
# - Year 1 changed from 2002 to 2012 everywhere it is referenced.
# - The "Average Values"/"Riddership Effect" source data on Sheet1 was refreshed
#   with new model output (new E/H column values for rows 8-18, 20, 21).
# - The % Diff formulas (columns G and I) switched from a *100 "percent-as-number"
#   convention to a native percentage (cell formatted as 0.00%), so the *100 is
#   dropped from every IFERROR(...,0) formula and the cells get a percent number
#   format instead.
# - The raw Average-Value cells (E/F columns) and the Total Modeled/Observed
#   Ridership rows got a #,##0.00 number format.
# - Selection/scroll position on Sheet1 moved.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Year label: 2002 -> 2012
# ---------------------------------------------------------------------------
$ws1.Range("C1").NumberFormat = "@"
$ws1.Range("C1").Value = "2012"

$ws1.Range("E7").NumberFormat = "@"
$ws1.Range("E7").Value = "2012"

# ---------------------------------------------------------------------------
# Factor rows 8-18: refreshed Average Values (E, H) + new % Diff formulas
# ---------------------------------------------------------------------------
$numFmt = "#,##0.00"
$pctFmt = "0.00%"

$factorRows = @{
    8  = @{ E = 4797772.19;          H = -929640.6071199999 }
    9  = @{ E = 1.11724546;          H = 134995.81056 }
    10 = @{ E = 1370736.25;          H = 135272.699867 }
    11 = @{ E = 34.76192261;         H = 33488.644388 }
    12 = @{ E = 4.1311;              H = -520223.48983 }
    13 = @{ E = 38883.57;            H = -11589.15046600002 }
    14 = @{ E = 8.5;                 H = -72121.204453 }
    15 = @{ E = 5.2;                 H = -21640.659565 }
    16 = @{ E = $null;               H = -1297488.2743 }
    17 = @{ E = 0;                   H = 0 }
    18 = @{ E = 0;                   H = 0 }
}

foreach ($r in 8..18) {
    $vals = $factorRows[$r]

    $eCell = $ws1.Range("E$r")
    $fCell = $ws1.Range("F$r")
    $eCell.NumberFormat = $numFmt
    $fCell.NumberFormat = $numFmt
    if ($null -ne $vals.E) {
        $eCell.Value = $vals.E
    }

    $gCell = $ws1.Range("G$r")
    $gCell.Formula = "=IFERROR((F$r-E$r)/E$r,0)"
    $gCell.NumberFormat = $pctFmt

    $hCell = $ws1.Range("H$r")
    $hCell.NumberFormat = $numFmt
    $hCell.Value = $vals.H

    $iCell = $ws1.Range("I$r")
    $iCell.Formula = "=IFERROR(H$r/`$E`$21,0)"
    $iCell.NumberFormat = $pctFmt
}

# ---------------------------------------------------------------------------
# Row 19 "New Reporters": D/E/F get the new number format; G/I formula+format;
# H19 now carries an explicit 0 value.
# ---------------------------------------------------------------------------
$ws1.Range("D19").NumberFormat = $numFmt
$ws1.Range("E19").NumberFormat = $numFmt
$ws1.Range("F19").NumberFormat = $numFmt

$ws1.Range("G19").Formula = "=IFERROR((F19-E19)/E19,0)"
$ws1.Range("G19").NumberFormat = $pctFmt

$ws1.Range("H19").Value = 0

$ws1.Range("I19").Formula = "=IFERROR(H19/`$E`$21,0)"
$ws1.Range("I19").NumberFormat = $pctFmt

# ---------------------------------------------------------------------------
# Row 20 "Total Modeled Ridership": refreshed E value, new % Diff format.
# ---------------------------------------------------------------------------
$ws1.Range("E20").NumberFormat = $numFmt
$ws1.Range("E20").Value = 10266082.9
$ws1.Range("F20").NumberFormat = $numFmt
$ws1.Range("H20").NumberFormat = $numFmt

$ws1.Range("G20").Formula = "=IFERROR((F20-E20)/E20,0)"
$ws1.Range("G20").NumberFormat = $pctFmt

$ws1.Range("I20").Formula = "=G20"
$ws1.Range("I20").NumberFormat = $pctFmt

# ---------------------------------------------------------------------------
# Row 21 "Total Observed Ridership": refreshed E value, new % Diff format.
# ---------------------------------------------------------------------------
$ws1.Range("E21").NumberFormat = $numFmt
$ws1.Range("E21").Value = 10128327.81
$ws1.Range("F21").NumberFormat = $numFmt
$ws1.Range("H21").NumberFormat = $numFmt

$ws1.Range("G21").Formula = "=IFERROR((F21-E21)/E21,0)"
$ws1.Range("G21").NumberFormat = $pctFmt

$ws1.Range("I21").Formula = "=G21"
$ws1.Range("I21").NumberFormat = $pctFmt

# ---------------------------------------------------------------------------
# Selection / scroll position on Sheet1 moved from K20 to H21 and the frozen
# top-left scroll cell (A7) was reset back to the default.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null
$ws1.Range("H21").Select() | Out-Null

Write-Host "Sheet1 factor table refreshed for Year 1 = 2012."
